# Edit script: v3.0 update FCI 27/1/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add date header in column C (B1 stays "06-01-2023")
# Copy B1's format (bold, centered, thin border) onto C1, then set its value
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = "13-01-2023"

# Rows 2-52: labels (column A)
$labels = @(
    "1810 Renta variable",
    "1822 Raices Valores Negociables",
    "Adcap IOL Acciones Argentina",
    "Allaria Acciones",
    "Alpha Acciones",
    "Alpha Mega",
    "Alpha Mercosur",
    "Alpha Recursos Naturales",
    "Alpha planeam equil",
    "Alpha renta balan global",
    "Argenfunds",
    "Arpenta acciones",
    "Arpenta ex Mercosur",
    "Balanz",
    "Bull Market",
    "CMA acciones",
    "Compass Crecimiento",
    "Consultatio Acciones Argentina",
    "Consultatio Renta Variable",
    "Delta Acciones",
    "Delta Recursos Naturales",
    "Delta Select",
    "Delta gestion V",
    "FBA Acciones Argentinas",
    "FBA Calificado",
    "Fima Acciones",
    "Fima PB Acciones",
    "Gainvest Renta Variable",
    "Galileo Acciones",
    "Goal Acciones Argentinas",
    "Goal acciones plus",
    "HF Acciones Argentinas",
    "HF Acciones Lideres",
    "IAM Renta Variable",
    "Lombardi",
    "MAF",
    "Megainver",
    "Pellegrini Acciones",
    "Pionero Acciones",
    "Premier Renta Variable",
    "Quinquela Acciones",
    "Rofex 20 Renta Variable",
    "SBS Acciones Argentina",
    "Schroeder RV",
    "Supefondo RV",
    "Superfondo ",
    "Supergestion",
    "Toronto Trust Multimercado",
    "Toronto trust Argy",
    "avg",
    "total"
)

# Rows 2-52: column B values (as of 06-01-2023)
$bvals = @(
    78476.92,
    164475.4,
    26608.73,
    34819.78,
    85824.98,
    194733.06,
    95157.3,
    57183.19,
    1270.11,
    281306.99,
    8835.94,
    451.23,
    3894.07,
    143087.38,
    18210.62,
    29100.57,
    460670.55,
    299078.34,
    104654.6,
    17825.12,
    85191.82,
    236715.38,
    35857.41,
    122389.53,
    120543.2,
    196570.7,
    128858.05,
    71116.75,
    540169.14,
    26072.86,
    4065.63,
    88625.83,
    137116.83,
    28550.44,
    22630.17,
    15721.9,
    25097.35,
    59798.1,
    98012.9,
    35588.62,
    80197.7,
    56674.56,
    223662.49,
    241929.54,
    603943.27,
    88367.8,
    30657.53,
    25084.87,
    10957.58,
    113180.26,
    5545832.83
)

# Rows 2-52: column C values (as of 13-01-2023)
$cvals = @(
    90812.03,
    173984.09,
    27374.89,
    35876.99,
    85803.85,
    194718.8,
    95138.52,
    61249.9,
    1270.69,
    281236.14,
    8836.12,
    450.83,
    3894.72,
    148848.48,
    19324.04,
    29065.99,
    455957.45,
    299099.85,
    104680.76,
    17840.26,
    85194.27,
    236685.76,
    35804.91,
    118545.95,
    116901.71,
    199818.62,
    121976.55,
    71086.22,
    550171.07,
    26066.84,
    4066.04,
    90634.21,
    139635.16,
    29876.95,
    23803.54,
    15723.67,
    25101.66,
    59785.34,
    98045.22,
    35601.34,
    80186.7,
    56675.38,
    223571.26,
    242001.9,
    660221.95,
    88426.83,
    30750.72,
    25868.73,
    10959.5,
    115074.54,
    5638652.4
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $acell = $ws.Cells.Item($row, 1)
    $acell.Value = $labels[$i]
    $acell.Font.Bold = $true
    $acell.HorizontalAlignment = -4108
    $acell.VerticalAlignment = -4160
    $acell.Borders.LineStyle = 1
    $acell.Borders.Weight = 2
    $ws.Cells.Item($row, 2).Value = $bvals[$i]
    $ws.Cells.Item($row, 3).Value = $cvals[$i]
}

Write-Host "Edit applied successfully"
